# Report regenerated for archive:
#  - Status changes from "Ready for handoff" to "In Translation" for every
#    localized file row (Overview sheet's per-language status columns, plus
#    the Status column on each per-language sheet).
#  - The Status column is narrowed on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text ("Ready for handoff" -> "In Translation") ---
$overview.Range("E2:F3").Value = "In Translation"
$zhcn.Range("C2:C3").Value = "In Translation"
$dede.Range("C2:C3").Value = "In Translation"

# --- Narrow the Status columns ---
# target stored column width (OOXML <col width>) is 13.4101845877511;
# the ColumnWidth COM property below is the closest input value that the
# engine maps to that stored width.
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
